# Gym Workouts update: append 21 new workout log rows (rows 322-342)
# for exercise sessions on 2018-01-07 (Sunday) and 2018-01-08 (Monday).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data: ExerciseId, DateId, ExerciseDate(serial), ExerciseWeek, ExerciseMonth,
#       ExerciseYear, ExerciseDay, ExerciseName, Weight, Sets, Reps, TrainingArea
$rows = @(
    @(321, 38, 43107, 1, "January", 2018, "Sunday", "Leg Extension", 100, 4, 12, "Legs"),
    @(322, 38, 43107, 1, "January", 2018, "Sunday", "Tricep Pull down", 40, 4, 8, "Arms"),
    @(323, 38, 43107, 1, "January", 2018, "Sunday", "Leg Curl", 60, 4, 12, "Legs"),
    @(324, 38, 43107, 1, "January", 2018, "Sunday", "Hammer Curl", 20, 4, 8, "Arms"),
    @(325, 38, 43107, 1, "January", 2018, "Sunday", "Hip adduction", 45, 4, 12, "Legs"),
    @(326, 38, 43107, 1, "January", 2018, "Sunday", "Lying Bicep Curl", 30, 4, 8, "Arms"),
    @(327, 38, 43107, 1, "January", 2018, "Sunday", "Hip abduction", 52, 4, 12, "Legs"),
    @(328, 38, 43107, 1, "January", 2018, "Sunday", "Barbell twists", 0, 4, 12, "Core"),
    @(329, 38, 43107, 1, "January", 2018, "Sunday", "Left Crunch", 0, 4, 12, "Core"),
    @(330, 38, 43107, 1, "January", 2018, "Sunday", "Right Crunch", 0, 4, 12, "Core"),
    @(331, 38, 43107, 1, "January", 2018, "Sunday", "Leg Raises", 0, 4, 10, "Core"),
    @(332, 38, 43107, 1, "January", 2018, "Sunday", "Sled Pushes", 25, 5, 2, "Core"),
    @(333, 39, 43108, 2, "January", 2018, "Monday", "Incline Bench", 80, 5, 5, "Chest"),
    @(334, 39, 43108, 2, "January", 2018, "Monday", "Seated Row", 66, 4, 8, "Back"),
    @(335, 39, 43108, 2, "January", 2018, "Monday", "Overhead Press", 40, 5, 5, "Chest"),
    @(336, 39, 43108, 2, "January", 2018, "Monday", "Sit ups", 0, 5, 12, "Core"),
    @(337, 39, 43108, 2, "January", 2018, "Monday", "Bicycles", 0, 4, 12, "Core"),
    @(338, 39, 43108, 2, "January", 2018, "Monday", "Left Crunch", 0, 4, 10, "Core"),
    @(339, 39, 43108, 2, "January", 2018, "Monday", "Right Crunch", 0, 4, 10, "Core"),
    @(340, 39, 43108, 2, "January", 2018, "Monday", "Heel-taps", 0, 4, 10, "Core"),
    @(341, 39, 43108, 2, "January", 2018, "Monday", "Plank", 0, 3, 1, "Core")

)

$startRow = 322
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
    $ws.Cells.Item($r, 8).Value = $data[7]
    $ws.Cells.Item($r, 9).Value = $data[8]
    $ws.Cells.Item($r, 10).Value = $data[9]
    $ws.Cells.Item($r, 11).Value = $data[10]
    $ws.Cells.Item($r, 12).Value = $data[11]
}

# Match the saved view state from the edit: selection sits a few rows below
# the last data row (as left by the person entering the data).
$ws.Range("C345").Select()

Write-Output "Added $($rows.Count) rows ($startRow..$($startRow + $rows.Count - 1))"
